# zdemo_excel15#1 - regenerate the "generated at" timestamp cells, drop the
# explicit zero row heights left over from the previous export, and turn the
# "Click here to visit abap2xlsx homepage" label into a real hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The demo stamps a date (B3) / time (C3) pair when it is generated; refresh
# them to the new snapshot values.
$ws.Range("B3").Value = 44534
$ws.Range("C3").Value = 0.4749884259259259

# Every row in this sheet was being forced to height 0 by the previous
# exporter. Auto-fitting restores the normal/default row height, which drops
# the explicit ht="0" override from each <row> element.
$ws.UsedRange.Rows.AutoFit()

# Turn the existing "Click here to visit abap2xlsx homepage" text (B4) into a
# clickable hyperlink pointing at the project's homepage.
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/abap2xlsx/abap2xlsx")

# Adding the hyperlink auto-applies Excel's built-in blue/underline
# "Hyperlink" style; the cell keeps its original (unstyled) formatting, so
# put it back.
$ws.Range("B4").Style = "Normal"
